$wb = $excel.ActiveWorkbook

$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# Insert a new "State" column into hotel_info between Hotel_Name (B) and City (C)
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# Move review_info sheet so it comes before hotel_info (swap tab order)
$reviewSheet.Move($hotelSheet)

$wb.Save()
